$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-PlainValue {
    param($cellRef, $value)
    $ws.Range($cellRef).Value = $value
}

# Row 2
Set-PlainValue 'D2' '29.264.64'
Set-PlainValue 'E2' '  -0.28%  '

# Row 3
Set-PlainValue 'D3' '1.829.04'
Set-PlainValue 'E3' '  -0.58%  '

# Row 4
Set-TextValue 'D4' '1.006'
Set-PlainValue 'E4' '  +0.56%  '

# Row 5
Set-TextValue 'D5' '234.34'
Set-PlainValue 'E5' '  -1.97%  '

# Row 6
Set-TextValue 'D6' '0.5992'
Set-PlainValue 'E6' '  -4.22%  '

# Row 7
Set-PlainValue 'E7' '  +0.39%  '

# Row 8
Set-TextValue 'D8' '0.06969'
Set-PlainValue 'E8' '  -5.59%  '

# Row 9
Set-TextValue 'D9' '0.2757'
Set-PlainValue 'E9' '  -4.42%  '

# Row 10
Set-TextValue 'D10' '23.27'
Set-PlainValue 'E10' '  -5.99%  '

# Row 12
Set-PlainValue 'D12' '1.831.55'
Set-PlainValue 'E12' '  -0.38%  '

# Row 13
Set-TextValue 'D13' '4.749'
Set-PlainValue 'E13' '  -4.19%  '

# Row 14
Set-TextValue 'D14' '0.6258'
Set-PlainValue 'E14' '  -6.08%  '

# Row 15
Set-TextValue 'D15' '0.000009623'
Set-PlainValue 'E15' '  -7.12%  '

# Row 16
Set-TextValue 'D16' '78.40'
Set-PlainValue 'E16' '  -3.73%  '

# Row 17
Set-PlainValue 'D17' '28.513.34'
Set-PlainValue 'E17' '  -2.76%  '

# Row 18
Set-TextValue 'D18' '5.603'
Set-PlainValue 'E18' '  -10.12%  '

# Row 19
Set-TextValue 'D19' '219.75'
Set-PlainValue 'E19' '  -6.20%  '

# Row 20
Set-TextValue 'D20' '1.005'
Set-PlainValue 'E20' '  +0.39%  '

# Row 21
Set-TextValue 'D21' '11.55'
Set-PlainValue 'E21' '  -5.84%  '

# Row 22
Set-TextValue 'D22' '6.862'
Set-PlainValue 'E22' '  -5.80%  '

# Row 23
Set-TextValue 'D23' '1.007'
Set-PlainValue 'E23' '  +0.56%  '

# Row 24
Set-TextValue 'D24' '156.43'

# Row 25
Set-TextValue 'D25' '7.949'
Set-PlainValue 'E25' '  -5.96%  '

# Row 26
Set-TextValue 'D26' '0.1285'
Set-PlainValue 'E26' '  -3.74%  '

# Row 27
Set-TextValue 'D27' '16.51'
Set-PlainValue 'E27' '  -4.43%  '

# Row 28
Set-TextValue 'D28' '1.453'
Set-PlainValue 'E28' '  -2.40%  '

# Row 29
Set-TextValue 'D29' '0.06333'
Set-PlainValue 'E29' '  -12.11%  '

# Row 30
Set-TextValue 'D30' '1.438'
Set-PlainValue 'E30' '  -2.90%  '

# Row 31
Set-TextValue 'D31' '3.833'
Set-PlainValue 'E31' '  -4.70%  '

# Row 32
Set-TextValue 'D32' '3.748'
Set-PlainValue 'E32' '  -6.94%  '

# Row 33
Set-PlainValue 'E33' '  -5.00%  '

# Row 34
Set-TextValue 'D34' '1.087'
Set-PlainValue 'E34' '  -5.73%  '

# Row 35
Set-TextValue 'D35' '0.6437'
Set-PlainValue 'E35' '  -8.79%  '

# Row 36
Set-TextValue 'D36' '2.544'
Set-PlainValue 'E36' '  -1.66%  '

# Row 37
Set-TextValue 'D37' '2.753'
Set-PlainValue 'E37' '  -1.12%  '

# Row 38
Set-TextValue 'D38' '0.01757'
Set-PlainValue 'E38' '  -4.04%  '

# Row 39
Set-TextValue 'D39' '6.580'
Set-PlainValue 'E39' '  -2.61%  '

# Row 40
Set-PlainValue 'D40' '1.153.31'
Set-PlainValue 'E40' '  -6.39%  '

# Row 41
Set-TextValue 'D41' '0.8904'

# Row 42
Set-TextValue 'D42' '1.006'
Set-PlainValue 'E42' '  +0.44%  '

# Row 43
Set-PlainValue 'D43' '1.985.58'
Set-PlainValue 'E43' '  +0.16%  '

# Row 44
Set-TextValue 'D44' '100.49'
Set-PlainValue 'E44' '  -0.49%  '

# Row 45
Set-TextValue 'D45' '61.94'
Set-PlainValue 'E45' '  -4.75%  '

# Row 46
Set-PlainValue 'E46' '  -3.63%  '

# Row 47
Set-TextValue 'D47' '1.590'
Set-PlainValue 'E47' '  -5.82%  '

# Row 48
Set-PlainValue 'B48' 'EnergySwap'
Set-PlainValue 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '8.469'
Set-PlainValue 'E48' '  -4.45%  '

# Row 49
Set-PlainValue 'B49' 'Cronos'
Set-PlainValue 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D49' '0.05549'
Set-PlainValue 'E49' '  -1.67%  '

# Row 50
Set-TextValue 'D50' '0.4563'
Set-PlainValue 'E50' '  -0.26%  '

# Row 51
Set-TextValue 'D51' '6.403'
Set-PlainValue 'E51' '  -7.56%  '
